$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume/issue number and week-covering date range --------------
$ws.Range("A8").Value = "Volume 31   Number  10"
$ws.Range("C9").Value = "Report Covering the Week  3/4/2024  Through  3/10/2024"

# --- Crime Complaints table (rows 15-30): refreshed weekly figures ----------
# Row 15 (Rape)
$ws.Range("C15:E15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("C14:E14").Copy()
$ws.Range("C15:E15").PasteSpecial(-4122)
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 2
$ws.Range("J15").Value = 2
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -33.333333333333
$ws.Range("M15").Value = 100

# Row 16 (Robbery)
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -16.666666666666
$ws.Range("I16").Value = 27
$ws.Range("J16").Value = 33
$ws.Range("K16").Value = -18.181818181818
$ws.Range("L16").Value = -3.571428571428
$ws.Range("M16").Value = -27.027027027027

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -28.571428571428
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = -22.727272727272
$ws.Range("I17").Value = 51
$ws.Range("J17").Value = 49
$ws.Range("K17").Value = 4.081632653061
$ws.Range("L17").Value = 4.081632653061
$ws.Range("M17").Value = 64.516129032258

# Row 18 (Burglary)
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -54.166666666666
$ws.Range("I18").Value = 26
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = -35
$ws.Range("L18").Value = 30
$ws.Range("M18").Value = 44.444444444444

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 22
$ws.Range("H19").Value = 40.90909090909
$ws.Range("I19").Value = 63
$ws.Range("J19").Value = 53
$ws.Range("K19").Value = 18.867924528301
$ws.Range("L19").Value = 14.545454545454
$ws.Range("M19").Value = 61.538461538461

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -25
$ws.Range("I20").Value = 18
$ws.Range("J20").Value = 27
$ws.Range("K20").Value = -33.333333333333
$ws.Range("L20").Value = -45.454545454545
$ws.Range("M20").Value = 50

# Row 21 (TOTAL)
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 77
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = -15.384615384615
$ws.Range("I21").Value = 187
$ws.Range("J21").Value = 205
$ws.Range("K21").Value = -8.780487804878
$ws.Range("L21").Value = -1.058201058201
$ws.Range("M21").Value = 33.571428571428

# Row 22 (Transit)
$ws.Range("G22").Value = 1
$ws.Range("L22").Value = -57.142857142857

# Row 23 (Housing)
$ws.Range("C23").Value = 1
$ws.Range("D23").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 6
$ws.Range("J23").Value = 2
$ws.Range("K23").Value = 200
$ws.Range("L23").Value = -40
$ws.Range("M23").Value = 100

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 48
$ws.Range("G24").Value = 75
$ws.Range("H24").Value = -36
$ws.Range("I24").Value = 140
$ws.Range("J24").Value = 163
$ws.Range("K24").Value = -14.110429447852
$ws.Range("L24").Value = -50.17793594306
$ws.Range("M24").Value = 72.839506172839

# Row 25 (Retail Theft)
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -85.714285714285
$ws.Range("F25").Value = 5
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = -76.190476190476
$ws.Range("I25").Value = 16
$ws.Range("J25").Value = 44
$ws.Range("K25").Value = -63.636363636363
$ws.Range("L25").Value = -89.873417721519

# Row 26 (Misd. Assault)
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -37.5
$ws.Range("F26").Value = 22
$ws.Range("G26").Value = 20
$ws.Range("H26").Value = 10
$ws.Range("I26").Value = 60
$ws.Range("J26").Value = 70
$ws.Range("K26").Value = -14.285714285714
$ws.Range("L26").Value = -21.052631578947
$ws.Range("M26").Value = -24.050632911392

# Row 27 (UCR Rape*)
$ws.Range("C27:E27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("C14:E14").Copy()
$ws.Range("C27:E27").PasteSpecial(-4122)

# Row 28 (Other Sex Crimes)
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("E28").Value = "***.*"
$ws.Range("D14:E14").Copy()
$ws.Range("D28:E28").PasteSpecial(-4122)
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 200
$ws.Range("I28").Value = 9
$ws.Range("K28").Value = 12.5
$ws.Range("L28").Value = 0

# Row 29 (Shooting Vic.)
$ws.Range("D29").Value = 1
$ws.Range("G29").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = -100
$ws.Range("H29").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("G29").Value = 3
$ws.Range("J29").Value = 3

# Row 30 (Shooting Inc.)
$ws.Range("D30").Value = 1
$ws.Range("G30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("H30").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("G30").Value = 3
$ws.Range("J30").Value = 3
